# Apply the "update to config file to match dev branch config" change.
#
# Summary of content changes on the "Findings" worksheet (sheet5.xml):
#  1) Rows 39/40 swap their content (MID_MinorLessThan50Percent group moves
#     above MID_FindingError).
#  2) Five new rows are inserted, each directly above the "<Doc>_FindingError"
#     row for Paystub / VOE / Offer Letter / Work Number / Verification
#     Services, containing a new "<Doc>_MultipleIncomeDocsExists" check.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Findings")
$ws.Activate()

$multiDocDesc = "Multiple employment documents provided for same employment. Can't determine which employment document to use to validate Exact Day Calculator. Manual review required."

# --- 1) Swap rows 39 and 40 ---------------------------------------------
$a39 = $ws.Range("A39").Value()
$b39 = $ws.Range("B39").Value()
$c39 = $ws.Range("C39").Value()

$a40 = $ws.Range("A40").Value()
$b40 = $ws.Range("B40").Value()
$c40 = $ws.Range("C40").Value()

$ws.Range("A39").Value = $a40
$ws.Range("B39").Value = $b40
$ws.Range("C39").Value = $c40

$ws.Range("A40").Value = $a39
$ws.Range("B40").Value = $b39
$ws.Range("C40").Value = ""

# --- 2) Insert the five new "MultipleIncomeDocsExists" rows -------------
# Each row is inserted immediately above the existing "<Doc>_FindingError"
# row. Insert from the top down, so row numbers below are recalculated
# automatically by Excel as each insert shifts everything else down.

$ws.Range("A139").EntireRow.Insert()
$ws.Range("A139").Value = "PS_MultipleIncomeDocsExists"
$ws.Range("B139").Value = $multiDocDesc

$ws.Range("A161").EntireRow.Insert()
$ws.Range("A161").Value = "VOE_MultipleIncomeDocsExists"
$ws.Range("B161").Value = $multiDocDesc

$ws.Range("A176").EntireRow.Insert()
$ws.Range("A176").Value = "OL_MultipleIncomeDocsExists"
$ws.Range("B176").Value = $multiDocDesc

$ws.Range("A192").EntireRow.Insert()
$ws.Range("A192").Value = "WN_MultipleIncomeDocsExists"
$ws.Range("B192").Value = $multiDocDesc

$ws.Range("A206").EntireRow.Insert()
$ws.Range("A206").Value = "VS_MultipleIncomeDocsExists"
$ws.Range("B206").Value = $multiDocDesc

# --- 3) Update the view (scroll position / selection) -------------------
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 261
$excel.ActiveWindow.ScrollColumn = 1
